# edit.ps1 - apply "add email find password functionality" changes
#
# Summary of the edit (per the unified diff):
#   1. Prepend a "Front end " paragraph before the existing first
#      paragraph ("Angular: interceptor ...").
#   2. After the "Angular: ..." paragraph, insert a handful of new
#      notes ("Double direction bind", "Restfull update data", three
#      blank paragraphs, "backend").
#   3. At the end of the document (after the "... 验证失败" paragraph)
#      append a new numbered list of backend TODO items, finishing
#      with a blank list paragraph. The "_GoBack" bookmark - which
#      previously sat at the very end of the document - moves along
#      with the last typed text ("spring email implements reset
#      password by email"), matching Word's usual behaviour of
#      keeping _GoBack at the most recently edited spot.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------
# 1 & 2: rebuild the first paragraph plus its new neighbours.
# ---------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$headRange = $d.Range($firstPara.Range.Start, $firstPara.Range.End - 1)

$headXml = ""
$headXml += "<w:p $wNs><w:r><w:t xml:space='preserve'>Front end </w:t></w:r></w:p>"
$headXml += "<w:p $wNs><w:r><w:t>Angular: interceptor for add token to each request if url match the requirement</w:t></w:r></w:p>"
$headXml += "<w:p $wNs>"
$headXml +=   "<w:r><w:t xml:space='preserve'>Double direction </w:t></w:r>"
$headXml +=   "<w:proofErr w:type='gramStart'/>"
$headXml +=   "<w:r><w:t>bind</w:t></w:r>"
$headXml +=   "<w:proofErr w:type='gramEnd'/>"
$headXml +=   "<w:r><w:t xml:space='preserve'> </w:t></w:r>"
$headXml += "</w:p>"
$headXml += "<w:p $wNs><w:r><w:t>Restfull update data</w:t></w:r></w:p>"
$headXml += "<w:p $wNs/>"
$headXml += "<w:p $wNs/>"
$headXml += "<w:p $wNs/>"
$headXml += "<w:p $wNs><w:r><w:t>backend</w:t></w:r></w:p>"

$headRange.InsertXML($headXml)

# ---------------------------------------------------------------
# 3: append the new list items after the last paragraph, moving the
#    _GoBack bookmark onto the final typed item.
# ---------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$tailRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$tailXml = ""
# Unchanged "... 验证失败" paragraph, now without the bookmark (it
# moves to the new final item below).
$tailXml += "<w:p $wNs>"
$tailXml +=   "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>用</w:t></w:r>"
$tailXml +=   "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>throw</w:t></w:r>"
$tailXml +=   "<w:r><w:t xml:space='preserve'> </w:t></w:r>"
$tailXml +=   "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>excption</w:t></w:r>"
$tailXml +=   "<w:r><w:t xml:space='preserve'> to tell feign client, </w:t></w:r>"
$tailXml +=   "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>验证失败</w:t></w:r>"
$tailXml += "</w:p>"

$listPPr = "<w:pPr><w:pStyle w:val='a7'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>"

# "Restful  db passed"
$tailXml += "<w:p $wNs>$listPPr"
$tailXml +=   "<w:proofErr w:type='gramStart'/>"
$tailXml +=   "<w:r><w:t>Restful  db</w:t></w:r>"
$tailXml +=   "<w:proofErr w:type='gramEnd'/>"
$tailXml +=   "<w:r><w:t xml:space='preserve'> passed</w:t></w:r>"
$tailXml += "</w:p>"

# "Spring scheduler  check leave request"
$tailXml += "<w:p $wNs>$listPPr"
$tailXml +=   "<w:r><w:t xml:space='preserve'>Spring </w:t></w:r>"
$tailXml +=   "<w:proofErr w:type='gramStart'/>"
$tailXml +=   "<w:r><w:t>scheduler  check</w:t></w:r>"
$tailXml +=   "<w:proofErr w:type='gramEnd'/>"
$tailXml +=   "<w:r><w:t xml:space='preserve'> leave request</w:t></w:r>"
$tailXml += "</w:p>"

# "rabbitMq update admin leave approval info"
$tailXml += "<w:p $wNs>$listPPr"
$tailXml +=   "<w:r><w:t>rabbitMq update admin leave approval info</w:t></w:r>"
$tailXml += "</w:p>"

# "spring security authenticationo and authorization"
$tailXml += "<w:p $wNs>$listPPr"
$tailXml +=   "<w:r><w:t>spring security authenticationo and authorization</w:t></w:r>"
$tailXml += "</w:p>"

# "jwt "
$tailXml += "<w:p $wNs>$listPPr"
$tailXml +=   "<w:r><w:t xml:space='preserve'>jwt </w:t></w:r>"
$tailXml += "</w:p>"

# "pass http header in microserver by @requestheader   to avoid interceptor's thread safe problem"
$tailXml += "<w:p $wNs>$listPPr"
$tailXml +=   "<w:r><w:t>pass http header in microserver by @requestheader   to avoid interceptor’s thread safe problem</w:t></w:r>"
$tailXml += "</w:p>"

# "spring email implements reset password by email" - the _GoBack
# bookmark lands here, at the spot of the most recent edit.
$tailXml += "<w:p $wNs>$listPPr"
$tailXml +=   "<w:r><w:t>spring email implements reset password by email</w:t></w:r>"
$tailXml +=   "<w:bookmarkStart w:id='0' w:name='_GoBack'/>"
$tailXml +=   "<w:bookmarkEnd w:id='0'/>"
$tailXml += "</w:p>"

# trailing blank list paragraph
$tailXml += "<w:p $wNs><w:pPr><w:pStyle w:val='a7'/></w:pPr></w:p>"

$tailRange.InsertXML($tailXml)
